$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels - reordered
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "kitchens_1"
$ws.Range("D1").Value = "living_rooms_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "kitchens_2"

# Row 2: shift the 1 from D2 to E2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

# Row 3: move the 1 from A3 to F3
$ws.Range("A3").Value = 0
$ws.Range("F3").Value = 1

# Row 4: move the 1 from F4 to A4
$ws.Range("A4").Value = 1
$ws.Range("F4").Value = 0

# Row 7: move the 1 from E7 to D7
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
